$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16.39505822365226
$ws.Range("C2").Value = 8.525222701492536
$ws.Range("D2").Value = 5.390011855418391
$ws.Range("E2").Value = 11.57574284004636
$ws.Range("F2").Value = 49.57651527291279
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("J2").Value = 10.293943827033
$ws.Range("K2").Value = 15.97252694999251
$ws.Range("M2").Value = 18.23634088441713
$ws.Range("N2").Value = 24.59824231507136

$ws.Range("B3").Value = 16.23565866884602
$ws.Range("C3").Value = 8.425275859132425
$ws.Range("D3").Value = 5.396439464675177
$ws.Range("E3").Value = 11.58537427125258
$ws.Range("F3").Value = 49.48887899350181
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("J3").Value = 10.31127922875074
$ws.Range("K3").Value = 15.87701324052378
$ws.Range("M3").Value = 18.21572541783846
$ws.Range("N3").Value = 24.63127097962112

$ws.Range("B4").Value = 16.14159039517525
$ws.Range("C4").Value = 8.365915977353515
$ws.Range("D4").Value = 5.401318553319227
$ws.Range("E4").Value = 11.59314335261705
$ws.Range("F4").Value = 49.44419366543287
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("J4").Value = 10.3230180125713
$ws.Range("K4").Value = 15.82216627869775
$ws.Range("M4").Value = 18.206749774208
$ws.Range("N4").Value = 24.65339467434958

$ws.Range("B5").Value = 16.10425483877714
$ws.Range("C5").Value = 8.342260119897194
$ws.Range("D5").Value = 5.403541818846336
$ws.Range("E5").Value = 11.59677618763787
$ws.Range("F5").Value = 49.42828629178047
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("J5").Value = 10.32807719104952
$ws.Range("K5").Value = 15.80079021640747
$ws.Range("M5").Value = 18.20402125706596
$ws.Range("N5").Value = 24.66287388062705

$ws.Range("B6").Value = 16.09811674951884
$ws.Range("C6").Value = 8.338365141916729
$ws.Range("D6").Value = 5.403925197621175
$ws.Range("E6").Value = 11.59740762248293
$ws.Range("F6").Value = 49.42578413543087
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("J6").Value = 10.32893391189455
$ws.Range("K6").Value = 15.79730012768238
$ws.Range("M6").Value = 18.20362438668593
$ws.Range("N6").Value = 24.66447589548988

$ws.Range("B7").Value = 16.14108278051276
$ws.Range("C7").Value = 8.365594747687148
$ws.Range("D7").Value = 5.401347584913049
$ws.Range("E7").Value = 11.59319045559868
$ws.Range("F7").Value = 49.4439698021071
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("J7").Value = 10.32308512650855
$ws.Range("K7").Value = 15.82187402361765
$ws.Range("M7").Value = 18.20670921074441
$ws.Range("N7").Value = 24.65352063719063

$ws.Range("B8").Value = 16.33933335829807
$ws.Range("C8").Value = 8.49036196950145
$ws.Range("D8").Value = 5.392034872055294
$ws.Range("E8").Value = 11.57867890912977
$ws.Range("F8").Value = 49.54440998367838
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("J8").Value = 10.29969399727488
$ws.Range("K8").Value = 15.93881741590471
$ws.Range("M8").Value = 18.22847048921902
$ws.Range("N8").Value = 24.60924784947331

$ws.Range("B9").Value = 16.75636212660863
$ws.Range("C9").Value = 8.749640890341789
$ws.Range("D9").Value = 5.381147751617442
$ws.Range("E9").Value = 11.56492542711225
$ws.Range("F9").Value = 49.81334306824617
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("J9").Value = 10.26250021165269
$ws.Range("K9").Value = 16.19736501944239
$ws.Range("M9").Value = 18.30019734324548
$ws.Range("N9").Value = 24.53706314533306

$ws.Range("B10").Value = 17.07729917188322
$ws.Range("C10").Value = 8.947217428870372
$ws.Range("D10").Value = 5.37761047632593
$ws.Range("E10").Value = 11.56375586315977
$ws.Range("F10").Value = 50.05413698946126
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("J10").Value = 10.24044919917005
$ws.Range("K10").Value = 16.40380482387853
$ws.Range("M10").Value = 18.37034949442809
$ws.Range("N10").Value = 24.49295535631572

$ws.Range("B11").Value = 17.22586761958424
$ws.Range("C11").Value = 9.038242865189606
$ws.Range("D11").Value = 5.376961640402539
$ws.Range("E11").Value = 11.56515551078205
$ws.Range("F11").Value = 50.17288486241453
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("J11").Value = 10.2315600257243
$ws.Range("K11").Value = 16.50099017372081
$ws.Range("M11").Value = 18.40598589497839
$ws.Range("N11").Value = 24.47482904465743

$ws.Range("B12").Value = 17.28244610970791
$ws.Range("C12").Value = 9.072843642523711
$ws.Range("D12").Value = 5.376853225500333
$ws.Range("E12").Value = 11.56596233606292
$ws.Range("F12").Value = 50.21915784560949
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("J12").Value = 10.22835787249782
$ws.Range("K12").Value = 16.53823512940021
$ws.Range("M12").Value = 18.42000917043445
$ws.Range("N12").Value = 24.46824396038384

$ws.Range("B13").Value = 17.27024755273882
$ws.Range("C13").Value = 9.065386430792644
$ws.Range("D13").Value = 5.376870482185155
$ws.Range("E13").Value = 11.56577627740261
$ws.Range("F13").Value = 50.20913437115425
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("J13").Value = 10.22904022396393
$ws.Range("K13").Value = 16.53019451391107
$ws.Range("M13").Value = 18.41696561697131
$ws.Range("N13").Value = 24.46964976628754

$ws.Range("B14").Value = 17.23051626872394
$ws.Range("C14").Value = 9.041087048783776
$ws.Range("D14").Value = 5.376949974443409
$ws.Range("E14").Value = 11.56521634710809
$ws.Range("F14").Value = 50.17666571197282
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("J14").Value = 10.23129329749936
$ws.Range("K14").Value = 16.50404562373217
$ws.Range("M14").Value = 18.4071290473903
$ws.Range("N14").Value = 24.47428169398234

$ws.Range("B15").Value = 17.20621971554717
$ws.Range("C15").Value = 9.026219101866683
$ws.Range("D15").Value = 5.377016519534038
$ws.Range("E15").Value = 11.56490939148084
$ws.Range("F15").Value = 50.15694718669474
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("J15").Value = 10.23269471876099
$ws.Range("K15").Value = 16.48808551836724
$ws.Range("M15").Value = 18.40117248018339
$ws.Range("N15").Value = 24.47715521764894

$ws.Range("B16").Value = 17.06763727638195
$ws.Range("C16").Value = 8.941288910200464
$ws.Range("D16").Value = 5.3776721353727
$ws.Range("E16").Value = 11.56370318385771
$ws.Range("F16").Value = 50.04656030631158
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("J16").Value = 10.24105308724868
$ws.Range("K16").Value = 16.3975170611489
$ws.Range("M16").Value = 18.36809492574411
$ws.Range("N16").Value = 24.49417897862655

$ws.Range("B17").Value = 16.98324403330262
$ws.Range("C17").Value = 8.889456810174801
$ws.Range("D17").Value = 5.378319733818673
$ws.Range("E17").Value = 11.56345732070474
$ws.Range("F17").Value = 49.9811876296214
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("J17").Value = 10.24647299653725
$ws.Range("K17").Value = 16.34277540055674
$ws.Range("M17").Value = 18.34875193894439
$ws.Range("N17").Value = 24.50511915143101

$ws.Range("B18").Value = 16.93494849097582
$ws.Range("C18").Value = 8.859754300381658
$ws.Range("D18").Value = 5.378782653834426
$ws.Range("E18").Value = 11.56349779238326
$ws.Range("F18").Value = 49.94445505880815
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("J18").Value = 10.24969788419571
$ws.Range("K18").Value = 16.31159942055045
$ws.Range("M18").Value = 18.33797732429348
$ws.Range("N18").Value = 24.51159410257006

$ws.Range("B19").Value = 16.91864012132288
$ws.Range("C19").Value = 8.849717449153189
$ws.Range("D19").Value = 5.378954947042791
$ws.Range("E19").Value = 11.5635427617024
$ws.Range("F19").Value = 49.93216767038025
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("J19").Value = 10.25080824556216
$ws.Range("K19").Value = 16.30109784558191
$ws.Range("M19").Value = 18.33438970392435
$ws.Range("N19").Value = 24.5138177402519

$ws.Range("B20").Value = 16.99220282872921
$ws.Range("C20").Value = 8.894963290668381
$ws.Range("D20").Value = 5.378241441451797
$ws.Range("E20").Value = 11.56346467537573
$ws.Range("F20").Value = 49.98805695083216
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("J20").Value = 10.24588491381168
$ws.Range("K20").Value = 16.34857086239083
$ws.Range("M20").Value = 18.35077475829679
$ws.Range("N20").Value = 24.50393566684647

$ws.Range("B21").Value = 17.24217806401079
$ws.Range("C21").Value = 9.048221064890322
$ws.Range("D21").Value = 5.376922906440853
$ws.Range("E21").Value = 11.56537330802843
$ws.Range("F21").Value = 50.18616726315
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("J21").Value = 10.23062706656084
$ws.Range("K21").Value = 16.51171439117703
$ws.Range("M21").Value = 18.41000400010239
$ws.Range("N21").Value = 24.47291361216663

$ws.Range("B22").Value = 17.40738613245043
$ws.Range("C22").Value = 9.149134643911758
$ws.Range("D22").Value = 5.376860958288778
$ws.Range("E22").Value = 11.56823362521543
$ws.Range("F22").Value = 50.3232448729481
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("J22").Value = 10.22161089307629
$ws.Range("K22").Value = 16.62090660858227
$ws.Range("M22").Value = 18.45179084824575
$ws.Range("N22").Value = 24.45426496992797

$ws.Range("B23").Value = 17.31906046751379
$ws.Range("C23").Value = 9.095217441733574
$ws.Range("D23").Value = 5.376821120164934
$ws.Range("E23").Value = 11.56655979081354
$ws.Range("F23").Value = 50.24939504236443
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("J23").Value = 10.22633562400845
$ws.Range("K23").Value = 16.5624030822569
$ws.Range("M23").Value = 18.42920931345709
$ws.Range("N23").Value = 24.46406925136293

$ws.Range("B24").Value = 16.98815185801112
$ws.Range("C24").Value = 8.89247350784669
$ws.Range("D24").Value = 5.378276555152352
$ws.Range("E24").Value = 11.56346078391768
$ws.Range("F24").Value = 49.98494867834285
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("J24").Value = 10.24615044675654
$ws.Range("K24").Value = 16.34594981081836
$ws.Range("M24").Value = 18.34985916363458
$ws.Range("N24").Value = 24.5044701430168

$ws.Range("B25").Value = 16.6407914186623
$ws.Range("C25").Value = 8.678121090375967
$ws.Range("D25").Value = 5.383306808568723
$ws.Range("E25").Value = 11.56707445192789
$ws.Range("F25").Value = 49.7329463378708
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("J25").Value = 10.27163468487243
$ws.Range("K25").Value = 16.12442470455844
$ws.Range("M25").Value = 18.27770701164929
$ws.Range("N25").Value = 24.55502383133981

